# Apply the "GIT UPDATE" change described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The cell E8 previously held the shared string "Good Morning" (which is
# removed from the shared-strings table entirely). Its text is replaced
# with a brand new shared string "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Excel also records the newly selected/active cell in the sheet view.
$ws.Range("E8").Select()
